$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2285.9375
$ws.Range("I113").Value = 2523.75
$ws.Range("J113").Value = 2206.6667
$ws.Range("K113").Value = 2523.75
$ws.Range("L113").Value = 2206.6667
$ws.Range("M113").Value = 730.25
$ws.Range("N113").Value = -8714.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 7813.4546
$ws.Range("I37").Value = 1111.3334
$ws.Range("K37").Value = 1111.3334
$ws.Range("M37").Value = -838.3334

$ws.Range("H45").Value = 1049.2188
$ws.Range("I45").Value = 826.5789
$ws.Range("J45").Value = 1374.6154
$ws.Range("K45").Value = 826.5789
$ws.Range("L45").Value = 1374.6154
$ws.Range("M45").Value = -449.5789
$ws.Range("N45").Value = -2128.6154

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H132").Value = 1132108.6
$ws.Range("I132").Value = 1924381.4
$ws.Range("J132").Value = 102154
$ws.Range("K132").Value = 5773144.199999999
$ws.Range("L132").Value = 306462
$ws.Range("M132").Value = -5770614.199999999
$ws.Range("N132").Value = -311522

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1364.6897
$ws.Range("I86").Value = 1363.1
$ws.Range("J86").Value = 1368.2222
$ws.Range("K86").Value = 1363.1
$ws.Range("L86").Value = 1368.2222
$ws.Range("M86").Value = -240.0999999999999
$ws.Range("N86").Value = -3614.2222

$ws.Range("H89").Value = 1364.6897
$ws.Range("I89").Value = 1363.1
$ws.Range("J89").Value = 1368.2222
$ws.Range("K89").Value = 6815.5
$ws.Range("L89").Value = 6841.111
$ws.Range("M89").Value = -1199.5
$ws.Range("N89").Value = -18073.111

$ws.Range("H99").Value = 1215
$ws.Range("I99").Value = 1190.8182
$ws.Range("J99").Value = 1303.6666
$ws.Range("K99").Value = 1190.8182
$ws.Range("L99").Value = 1303.6666
$ws.Range("M99").Value = 307.1818000000001
$ws.Range("N99").Value = -4299.6666

$ws.Range("H134").Value = 37113.645
$ws.Range("I134").Value = 47237.875
$ws.Range("J134").Value = 2402
$ws.Range("K134").Value = 141713.625
$ws.Range("L134").Value = 7206
$ws.Range("M134").Value = -139178.625
$ws.Range("N134").Value = -12276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3143.4243
$ws.Range("I31").Value = 1599.091
$ws.Range("J31").Value = 3915.5908
$ws.Range("K31").Value = 1599.091
$ws.Range("L31").Value = 3915.5908
$ws.Range("M31").Value = -1304.091
$ws.Range("N31").Value = -4505.5908

$ws.Range("H34").Value = 3143.4243
$ws.Range("I34").Value = 1599.091
$ws.Range("J34").Value = 3915.5908
$ws.Range("K34").Value = 1599.091
$ws.Range("L34").Value = 3915.5908
$ws.Range("M34").Value = -1397.091
$ws.Range("N34").Value = -4319.5908

$ws.Range("H58").Value = 2784.6934
$ws.Range("I58").Value = 1222.2632
$ws.Range("J58").Value = 4389.3516
$ws.Range("K58").Value = 1222.2632
$ws.Range("L58").Value = 4389.3516
$ws.Range("M58").Value = -1019.2632
$ws.Range("N58").Value = -4795.3516

$ws.Range("H136").Value = 2784.6934
$ws.Range("I136").Value = 1222.2632
$ws.Range("J136").Value = 4389.3516
$ws.Range("K136").Value = 3666.7896
$ws.Range("L136").Value = 13168.0548
$ws.Range("M136").Value = -1116.7896
$ws.Range("N136").Value = -18268.0548

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 223001500
$ws.Range("J105").Value = 223001500
$ws.Range("L105").Value = 669004500
$ws.Range("N105").Value = -669009742

$ws.Range("H131").Value = 4569
$ws.Range("I131").Value = 20120
$ws.Range("J131").Value = 2774.6538
$ws.Range("K131").Value = 60360
$ws.Range("L131").Value = 8323.9614
$ws.Range("M131").Value = -55320
$ws.Range("N131").Value = -18403.9614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4834.963
$ws.Range("I70").Value = 4764.6665
$ws.Range("J70").Value = 4891.2
$ws.Range("K70").Value = 4764.6665
$ws.Range("L70").Value = 4891.2
$ws.Range("M70").Value = -4494.6665
$ws.Range("N70").Value = -5431.2

$ws.Range("H73").Value = 4834.963
$ws.Range("I73").Value = 4764.6665
$ws.Range("J73").Value = 4891.2
$ws.Range("K73").Value = 4764.6665
$ws.Range("L73").Value = 4891.2
$ws.Range("M73").Value = -3828.6665
$ws.Range("N73").Value = -6763.2

$ws.Range("H126").Value = 1962
$ws.Range("I126").Value = 1646.909
$ws.Range("J126").Value = 2457.1428
$ws.Range("K126").Value = 4940.727000000001
$ws.Range("L126").Value = 7371.428400000001
$ws.Range("M126").Value = -2470.727000000001
$ws.Range("N126").Value = -12311.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3089.2144
$ws.Range("I16").Value = 3296.4614
$ws.Range("J16").Value = 395
$ws.Range("K16").Value = 3296.4614
$ws.Range("L16").Value = 395
$ws.Range("M16").Value = -3126.4614
$ws.Range("N16").Value = -735

$ws.Range("H22").Value = 414.2857
$ws.Range("I22").Value = 475
$ws.Range("J22").Value = 333.33334
$ws.Range("K22").Value = 475
$ws.Range("L22").Value = 333.33334
$ws.Range("M22").Value = -180
$ws.Range("N22").Value = -923.33334

$ws.Range("H27").Value = 414.2857
$ws.Range("I27").Value = 475
$ws.Range("J27").Value = 333.33334
$ws.Range("K27").Value = 475
$ws.Range("L27").Value = 333.33334
$ws.Range("M27").Value = -368
$ws.Range("N27").Value = -547.33334

$ws.Range("H40").Value = 1546.5555
$ws.Range("I40").Value = 1489.875
$ws.Range("K40").Value = 1489.875
$ws.Range("M40").Value = -1353.875

$ws.Range("H46").Value = 914.2105
$ws.Range("I46").Value = 858.7879
$ws.Range("K46").Value = 858.7879
$ws.Range("M46").Value = -670.7879

$ws.Range("H55").Value = 186.66667
$ws.Range("I55").Value = 180
$ws.Range("J55").Value = 200
$ws.Range("K55").Value = 180
$ws.Range("L55").Value = 200
$ws.Range("M55").Value = -7
$ws.Range("N55").Value = -546

$ws.Range("H136").Value = 2668.0833
$ws.Range("I136").Value = 1318.6666
$ws.Range("J136").Value = 4017.5
$ws.Range("K136").Value = 3955.9998
$ws.Range("L136").Value = 12052.5
$ws.Range("M136").Value = -1405.9998
$ws.Range("N136").Value = -17152.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 621.375
$ws.Range("I107").Value = 399.05554
$ws.Range("J107").Value = 1288.3334
$ws.Range("K107").Value = 1197.16662
$ws.Range("L107").Value = 3865.0002
$ws.Range("M107").Value = 722.83338
$ws.Range("N107").Value = -7705.0002

$ws.Range("H113").Value = 330.7647
$ws.Range("I113").Value = 328.2
$ws.Range("J113").Value = 350
$ws.Range("K113").Value = 984.5999999999999
$ws.Range("L113").Value = 1050
$ws.Range("M113").Value = 1185.4
$ws.Range("N113").Value = -5390

$ws.Range("H136").Value = 2426394.8
$ws.Range("I136").Value = 6548.095
$ws.Range("J136").Value = 8404839
$ws.Range("K136").Value = 19644.285
$ws.Range("L136").Value = 25214517
$ws.Range("M136").Value = -17094.285
$ws.Range("N136").Value = -25219617
